$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.153.12"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.08"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.77"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6577"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07392"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2919"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.84"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07738"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.842.10"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6649"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.73"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008425"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.158.99"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.093.13"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.06"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.121"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.585"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1392"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.91"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.516"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.110"
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.041"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.864"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7381"
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.141"
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.653"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.298.35"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.731"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9186"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("B42").Value = "XinFinNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.08596"
$ws.Range("E42").Value = "  +7.07%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.957"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.11"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.987.27"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5136"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.14"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05844"
$ws.Range("E51").Value = "  -1.28%  "
